# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 08:50"

# --- Finlandia / Grecia swap (row 42 / 43), with updated Finlandia stats ---
$ws.Cells.Item(42,1).Value = "Finlandia"
$ws.Cells.Item(42,2).Value = 1615
$ws.Cells.Item(42,3).Value = 97
$ws.Cells.Item(42,4).Value = 300
$ws.Cells.Item(42,5).Value = 1296
$ws.Cells.Item(42,6).Value = 65
$ws.Cells.Item(42,7).Value = 0
$ws.Cells.Item(42,8).Value = 19

$ws.Cells.Item(43,1).Value = "Grecia"
$ws.Cells.Item(43,2).Value = 1544
$ws.Cells.Item(43,3).Value = 0
$ws.Cells.Item(43,4).Value = 61
$ws.Cells.Item(43,5).Value = 1430
$ws.Cells.Item(43,6).Value = 91
$ws.Cells.Item(43,7).Value = 0
$ws.Cells.Item(43,8).Value = 53

# --- Banglades overtakes Aruba / Monaco / Madagascar (rows 128-131) ---
$ws.Cells.Item(128,1).Value = "Banglades"
$ws.Cells.Item(128,2).Value = 61
$ws.Cells.Item(128,3).Value = 5
$ws.Cells.Item(128,4).Value = 26
$ws.Cells.Item(128,5).Value = 29
$ws.Cells.Item(128,6).Value = 1
$ws.Cells.Item(128,7).Value = 0
$ws.Cells.Item(128,8).Value = 6

$ws.Cells.Item(129,1).Value = "Aruba"
$ws.Cells.Item(129,2).Value = 60
$ws.Cells.Item(129,3).Value = 0
$ws.Cells.Item(129,4).Value = 1
$ws.Cells.Item(129,5).Value = 59
$ws.Cells.Item(129,6).Value = 0
$ws.Cells.Item(129,7).Value = 0
$ws.Cells.Item(129,8).Value = 0

$ws.Cells.Item(130,1).Value = "Monaco"
$ws.Cells.Item(130,2).Value = 60
$ws.Cells.Item(130,3).Value = 0
$ws.Cells.Item(130,4).Value = 2
$ws.Cells.Item(130,5).Value = 57
$ws.Cells.Item(130,6).Value = 2
$ws.Cells.Item(130,7).Value = 0
$ws.Cells.Item(130,8).Value = 1

$ws.Cells.Item(131,1).Value = "Madagascar"
$ws.Cells.Item(131,2).Value = 59
$ws.Cells.Item(131,3).Value = 0
$ws.Cells.Item(131,4).Value = 0
$ws.Cells.Item(131,5).Value = 59
$ws.Cells.Item(131,6).Value = 6
$ws.Cells.Item(131,7).Value = 0
$ws.Cells.Item(131,8).Value = 0

# --- Simple numeric updates ---
# Row 25: Chequia
$ws.Cells.Item(25,2).Value = 3869
$ws.Cells.Item(25,3).Value = 11
$ws.Cells.Item(25,4).Value = 71
$ws.Cells.Item(25,5).Value = 3752
$ws.Cells.Item(25,6).Value = 77
$ws.Cells.Item(25,7).Value = 2
$ws.Cells.Item(25,8).Value = 46

# Row 33: Rumania
$ws.Cells.Item(33,5).Value = 2355
$ws.Cells.Item(33,7).Value = 1
$ws.Cells.Item(33,8).Value = 116

# Row 82: Republica de Macedonia
$ws.Cells.Item(82,6).Value = 8

# Row 108: Sri Lanka
$ws.Cells.Item(108,4).Value = 22
$ws.Cells.Item(108,5).Value = 125

# Row 109: Georgia
$ws.Cells.Item(109,4).Value = 27
$ws.Cells.Item(109,5).Value = 121

# Row 198: Butan
$ws.Cells.Item(198,4).Value = 2
$ws.Cells.Item(198,5).Value = 3
